# The project's screenshots are no longer needed - remove every inline
# picture from the document. Each picture is the sole run inside its own
# paragraph, so deleting the InlineShape removes just that run and leaves
# the (now empty) paragraph mark in place, matching the target edit.

$d = $word.ActiveDocument

while ($d.InlineShapes.Count -gt 0) {
    $d.InlineShapes.Item(1).Delete()
}

Write-Output ("Remaining InlineShapes: " + $d.InlineShapes.Count)
